$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: NewLoanInput
# ---------------------------------------------------------------------------
$wsLoanInput = $wb.Worksheets.Item("NewLoanInput")
[void]$wsLoanInput.Activate()
[void]$wsLoanInput.Range("B2").Select()

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
[void]$wsSummary.Activate()
$wsSummary.Range("A4").Value = 150
$wsSummary.Range("B4").Value = 50
[void]$wsSummary.Range("D4").Select()

# ---------------------------------------------------------------------------
# Sheet: Repayment schedule
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
[void]$wsRepay.Activate()

# Clear column O data cells (O2:O15) back to default (no content, no format)
[void]$wsRepay.Range("O2:O15").Clear()

# D3/E3 lose their special number formats (date / italic) and revert to the
# plain body style used elsewhere on the sheet (copy format from A3), then
# D3's leftover date value is cleared.
[void]$wsRepay.Range("A3").Copy()
[void]$wsRepay.Range("D3:E3").PasteSpecial(-4122)
[void]$wsRepay.Range("D3").ClearContents()
$excel.CutCopyMode = 0

# P3 keeps its style but its value is cleared
[void]$wsRepay.Range("P3").ClearContents()

# Value corrections
$wsRepay.Range("I3").Value = 50
$wsRepay.Range("L3").Value = 937.72

$wsRepay.Range("I5").Value = 0
$wsRepay.Range("K5").Value = 887.72
$wsRepay.Range("P5").Value = 887.72

$wsRepay.Range("I6").Value = 100
$wsRepay.Range("K6").Value = 987.72
$wsRepay.Range("P6").Value = 987.72

[void]$wsRepay.Range("P3").Select()

# ---------------------------------------------------------------------------
# Sheet: Transactions
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
[void]$wsTrans.Activate()

$wsTrans.Range("A2").Value = 6350
[void]$wsTrans.Range("J4").Copy()
[void]$wsTrans.Range("J2").PasteSpecial(-4122)
$wsTrans.Range("J2").Value = 10000

$wsTrans.Range("A3").Value = 6354
$wsTrans.Range("E3").Value = 937.72
$wsTrans.Range("F3").Value = 786.76
[void]$wsTrans.Range("F2").Copy()
[void]$wsTrans.Range("J3").PasteSpecial(-4122)
$wsTrans.Range("J3").Value = 0
$excel.CutCopyMode = 0

$wsTrans.Range("A4").Value = 691

[void]$wsTrans.Range("D4").Select()
